$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 8..15 (extr1..extr8 contingencies) down to rows
# 10..17, copying both value and format, bottom-up so we never clobber a
# row before it has been read. This preserves the original shared-string
# entries / cell styles for the "extr*" rows untouched.
for ($r = 15; $r -ge 8; $r--) {
    $srcRow = $r
    $dstRow = $r + 2
    $ws.Range("A$srcRow`:E$srcRow").Copy()
    $ws.Range("A$dstRow`:E$dstRow").PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# The copy-down above doesn't carry formatting onto rows beyond the
# original used range (16, 17), so reapply the header-style (same as
# every other column-A cell) explicitly there.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 8: line7 contingency (uses row 7's formatting as a template)
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(8, 1).Value2 = 6
$ws.Cells.Item(8, 2).Value2 = "line7"
$ws.Cells.Item(8, 3).Value2 = 14
$ws.Cells.Item(8, 4).Value2 = 11
$ws.Cells.Item(8, 5).Value2 = $true

# New row 9: line8 contingency
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(9, 1).Value2 = 7
$ws.Cells.Item(9, 2).Value2 = "line8"
$ws.Cells.Item(9, 3).Value2 = 16
$ws.Cells.Item(9, 4).Value2 = 9
$ws.Cells.Item(9, 5).Value2 = $true

# Update A-column sequence numbers and the C/D/E values for the shifted
# extr1..extr8 rows per the target state.
$ws.Cells.Item(10, 1).Value2 = 8
$ws.Cells.Item(10, 3).Value2 = 5
$ws.Cells.Item(10, 4).Value2 = 12
$ws.Cells.Item(10, 5).Value2 = $true

$ws.Cells.Item(11, 1).Value2 = 9
$ws.Cells.Item(11, 3).Value2 = 5
$ws.Cells.Item(11, 4).Value2 = 9
$ws.Cells.Item(11, 5).Value2 = $true

$ws.Cells.Item(12, 1).Value2 = 10
$ws.Cells.Item(12, 3).Value2 = 10
$ws.Cells.Item(12, 4).Value2 = 11
$ws.Cells.Item(12, 5).Value2 = $true

$ws.Cells.Item(13, 1).Value2 = 11
$ws.Cells.Item(13, 3).Value2 = 7
$ws.Cells.Item(13, 4).Value2 = 8
$ws.Cells.Item(13, 5).Value2 = $false

$ws.Cells.Item(14, 1).Value2 = 12
$ws.Cells.Item(14, 3).Value2 = 9
$ws.Cells.Item(14, 4).Value2 = 11
$ws.Cells.Item(14, 5).Value2 = $false

$ws.Cells.Item(15, 1).Value2 = 13
$ws.Cells.Item(15, 3).Value2 = 7
$ws.Cells.Item(15, 4).Value2 = 11
$ws.Cells.Item(15, 5).Value2 = $false

$ws.Cells.Item(16, 1).Value2 = 14
$ws.Cells.Item(16, 3).Value2 = 5
$ws.Cells.Item(16, 4).Value2 = 7
$ws.Cells.Item(16, 5).Value2 = $true

$ws.Cells.Item(17, 1).Value2 = 15
$ws.Cells.Item(17, 3).Value2 = 8
$ws.Cells.Item(17, 4).Value2 = 5
$ws.Cells.Item(17, 5).Value2 = $true
